$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reposition/resize the workbook window to match the recorded Excel session
$win = $wb.Windows.Item(1)
$win.Left = 5805
$win.Top = 1530
$win.Width = 21600
$win.Height = 11385

# Update F1, keep A1 as-is
$ws.Range("F1").Value = "Badland"

# Add new rows of data
$ws.Range("F2").Value = "Hello World"
$ws.Range("F3").Value = "Badland"

# Clear the old F4 value (3.14) - row becomes empty but row element stays
$ws.Range("F4").ClearContents()
$ws.Rows.Item(4).OutlineLevel = 0

# Add new value further down
$ws.Range("C6").Value = "Test1"

# Update selection to match target state
$ws.Range("G5").Select()
